$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 641.8182
$ws.Range("I32").Value = 545
$ws.Range("K32").Value = 545
$ws.Range("M32").Value = -219

$ws.Range("H38").Value = 1591.5
$ws.Range("J38").Value = 2962.5
$ws.Range("L38").Value = 8887.5
$ws.Range("N38").Value = -9631.5

$ws.Range("H39").Value = 1513.3334
$ws.Range("I39").Value = 74
$ws.Range("J39").Value = 2233
$ws.Range("K39").Value = 222
$ws.Range("L39").Value = 6699
$ws.Range("M39").Value = 74
$ws.Range("N39").Value = -7291

$ws.Range("H40").Value = 1214.9
$ws.Range("I40").Value = 1028.5714
$ws.Range("J40").Value = 1649.6666
$ws.Range("K40").Value = 1028.5714
$ws.Range("L40").Value = 1649.6666
$ws.Range("M40").Value = -853.5714
$ws.Range("N40").Value = -1999.6666

$ws.Range("H58").Value = 18903.834
$ws.Range("J58").Value = 21351.32
$ws.Range("L58").Value = 64053.96
$ws.Range("N58").Value = -64353.96

$ws.Range("H115").Value = 1530
$ws.Range("I115").Value = 795
$ws.Range("K115").Value = 2385
$ws.Range("M115").Value = -818

$ws.Range("H116").Value = 3494.3333
$ws.Range("I116").Value = 2942.457
$ws.Range("J116").Value = 4701.5625
$ws.Range("K116").Value = 2942.457
$ws.Range("L116").Value = 4701.5625
$ws.Range("M116").Value = 499.5430000000001
$ws.Range("N116").Value = -11585.5625

$ws.Range("H121").Value = 1160.1111
$ws.Range("J121").Value = 1913.25
$ws.Range("L121").Value = 5739.75
$ws.Range("N121").Value = -9233.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1925.36
$ws.Range("I32").Value = 1690.6957
$ws.Range("J32").Value = 4624
$ws.Range("K32").Value = 1690.6957
$ws.Range("L32").Value = 4624
$ws.Range("M32").Value = -1403.6957
$ws.Range("N32").Value = -5198

$ws.Range("H61").Value = 1620.38
$ws.Range("I61").Value = 763.69446
$ws.Range("J61").Value = 3823.2856
$ws.Range("K61").Value = 763.69446
$ws.Range("L61").Value = 3823.2856
$ws.Range("M61").Value = -551.69446
$ws.Range("N61").Value = -4247.2856

$ws.Range("H74").Value = 849.08
$ws.Range("I74").Value = 714.9091
$ws.Range("J74").Value = 1833
$ws.Range("K74").Value = 714.9091
$ws.Range("L74").Value = 1833
$ws.Range("M74").Value = 159.0909
$ws.Range("N74").Value = -3581

$ws.Range("H77").Value = 849.08
$ws.Range("I77").Value = 714.9091
$ws.Range("J77").Value = 1833
$ws.Range("K77").Value = 3574.5455
$ws.Range("L77").Value = 9165
$ws.Range("M77").Value = 793.4545000000003
$ws.Range("N77").Value = -17901

$ws.Range("H136").Value = 1620.38
$ws.Range("I136").Value = 763.69446
$ws.Range("J136").Value = 3823.2856
$ws.Range("K136").Value = 2291.08338
$ws.Range("L136").Value = 11469.8568
$ws.Range("M136").Value = 258.91662
$ws.Range("N136").Value = -16569.8568

$ws.Range("H138").Value = 49388.89
$ws.Range("J138").Value = 49388.89
$ws.Range("L138").Value = 49388.89
$ws.Range("N138").Value = -59668.89

$ws.Range("H139").Value = 26220.908
$ws.Range("J139").Value = 26220.908
$ws.Range("L139").Value = 26220.908
$ws.Range("N139").Value = -36500.908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1290.2858
$ws.Range("J52").Value = 1290.2858
$ws.Range("L52").Value = 3870.8574
$ws.Range("N52").Value = -4402.857400000001

$ws.Range("H59").Value = 2833.5
$ws.Range("I59").Value = 1101
$ws.Range("J59").Value = 3180
$ws.Range("K59").Value = 3303
$ws.Range("L59").Value = 9540
$ws.Range("M59").Value = -2763
$ws.Range("N59").Value = -10620

$ws.Range("H63").Value = 14533.333
$ws.Range("J63").Value = 15975
$ws.Range("L63").Value = 47925
$ws.Range("N63").Value = -49423

$ws.Range("H64").Value = 13066.777
$ws.Range("J64").Value = 19166.666
$ws.Range("L64").Value = 57499.99800000001
$ws.Range("N64").Value = -58039.99800000001

$ws.Range("H66").Value = 14533.333
$ws.Range("J66").Value = 15975
$ws.Range("L66").Value = 143775
$ws.Range("N66").Value = -151263

$ws.Range("H67").Value = 13066.777
$ws.Range("J67").Value = 19166.666
$ws.Range("L67").Value = 57499.99800000001
$ws.Range("N67").Value = -59371.99800000001

$ws.Range("H87").Value = 6327.727
$ws.Range("I87").Value = 2515
$ws.Range("J87").Value = 13000
$ws.Range("K87").Value = 7545
$ws.Range("L87").Value = 39000
$ws.Range("M87").Value = -6297
$ws.Range("N87").Value = -41496

$ws.Range("H90").Value = 6327.727
$ws.Range("I90").Value = 2515
$ws.Range("J90").Value = 13000
$ws.Range("K90").Value = 22635
$ws.Range("L90").Value = 117000
$ws.Range("M90").Value = -16395
$ws.Range("N90").Value = -129480

$ws.Range("H107").Value = 710.0244
$ws.Range("I107").Value = 258.3
$ws.Range("J107").Value = 855.74194
$ws.Range("K107").Value = 774.9000000000001
$ws.Range("L107").Value = 2567.22582
$ws.Range("M107").Value = 1145.1
$ws.Range("N107").Value = -6407.22582

$ws.Range("H116").Value = 1775
$ws.Range("I116").Value = 1033.3334
$ws.Range("K116").Value = 3100.0002
$ws.Range("M116").Value = 341.9998000000001

$ws.Range("H118").Value = 2370.9
$ws.Range("I118").Value = 996.3333
$ws.Range("K118").Value = 2988.9999
$ws.Range("M118").Value = -1745.9999

$ws.Range("H120").Value = 16635.889
$ws.Range("I120").Value = 12230
$ws.Range("K120").Value = 36690
$ws.Range("M120").Value = -31852

$ws.Range("H138").Value = 1738.6
$ws.Range("I138").Value = 902.6667
$ws.Range("K138").Value = 2708.0001
$ws.Range("M138").Value = 2431.9999

$ws.Range("H140").Value = 11908226
$ws.Range("I140").Value = 23810324
$ws.Range("J140").Value = 6128.5713
$ws.Range("K140").Value = 71430972
$ws.Range("L140").Value = 18385.7139
$ws.Range("M140").Value = -71425792
$ws.Range("N140").Value = -28745.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 29230.2
$ws.Range("J135").Value = 29230.2
$ws.Range("L135").Value = 29230.2
$ws.Range("N135").Value = -39370.2

$ws.Range("H141").Value = 37714.285
$ws.Range("J141").Value = 37714.285
$ws.Range("L141").Value = 37714.285
$ws.Range("N141").Value = -48074.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 51670
$ws.Range("J3").Value = 51670
$ws.Range("L3").Value = 51670
$ws.Range("N3").Value = -51894

$ws.Range("H7").Value = 1974.5454
$ws.Range("I7").Value = 1135.5555
$ws.Range("J7").Value = 5750
$ws.Range("K7").Value = 1135.5555
$ws.Range("L7").Value = 5750
$ws.Range("M7").Value = -1023.5555
$ws.Range("N7").Value = -5974

$ws.Range("H15").Value = 51670
$ws.Range("J15").Value = 51670
$ws.Range("L15").Value = 51670
$ws.Range("N15").Value = -52010

$ws.Range("H40").Value = 1919.5333
$ws.Range("I40").Value = 1171.6364
$ws.Range("J40").Value = 3976.25
$ws.Range("K40").Value = 1171.6364
$ws.Range("L40").Value = 3976.25
$ws.Range("M40").Value = -1035.6364
$ws.Range("N40").Value = -4248.25

$ws.Range("H46").Value = 2563.7273
$ws.Range("J46").Value = 3325
$ws.Range("L46").Value = 3325
$ws.Range("N46").Value = -3701

$ws.Range("H126").Value = 1974.5454
$ws.Range("I126").Value = 1135.5555
$ws.Range("J126").Value = 5750
$ws.Range("K126").Value = 3406.6665
$ws.Range("L126").Value = 17250
$ws.Range("M126").Value = -936.6664999999998
$ws.Range("N126").Value = -22190

$ws.Range("H132").Value = 2232.879
$ws.Range("I132").Value = 1440.909
$ws.Range("J132").Value = 3816.818
$ws.Range("K132").Value = 4322.727000000001
$ws.Range("L132").Value = 11450.454
$ws.Range("M132").Value = -1792.727000000001
$ws.Range("N132").Value = -16510.454

$ws.Range("H135").Value = 29764.047
$ws.Range("J135").Value = 29764.047
$ws.Range("L135").Value = 29764.047
$ws.Range("N135").Value = -39904.047

$ws.Range("H140").Value = 29259.818
$ws.Range("J140").Value = 29259.818
$ws.Range("L140").Value = 29259.818
$ws.Range("N140").Value = -39619.818

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4349491
$ws.Range("I126").Value = 914.2
$ws.Range("K126").Value = 2742.6
$ws.Range("M126").Value = -272.6000000000004

$ws.Range("H132").Value = 331605.3
$ws.Range("I132").Value = 504092.84
$ws.Range("J132").Value = 44126.082
$ws.Range("K132").Value = 1512278.52
$ws.Range("L132").Value = 132378.246
$ws.Range("M132").Value = -1509748.52
$ws.Range("N132").Value = -137438.246
